# Extended the email functionality and bug fixes:
#  - Rename "Email Distribution" header columns:
#       G1: "zip_file" -> "file"
#       H1: "Notes"    -> "Comments"
#  - Make "Email Distribution" the active sheet/tab again, with H1 selected
#    (previously "Email Settings" was the active tab with H5 selected).

$wb = $excel.ActiveWorkbook

$wsDist = $wb.Worksheets.Item("Email Distribution")
$wsDist.Range("G1").Value = "file"
$wsDist.Range("H1").Value = "Comments"

# Switch the active tab/selection back to the Email Distribution sheet.
$wsDist.Activate()
$wsDist.Range("H1").Select()
